$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 20 de Marzo de 2020 a las 22:46"

# Country column + stat columns (B:Casos totales C:Nuevos casos D:Casos activos
# E:Recuperados F:Casos criticos G:Muertes hoy H:Muertes) per affected row.
# Row 6: España -> España
$ws.Range("B6").Value = 21510
$ws.Range("C6").Value = 3433
$ws.Range("E6").Value = 18830
$ws.Range("G6").Value = 261
$ws.Range("H6").Value = 1092

# Row 7: Alemania -> Alemania
$ws.Range("E7").Value = 19600
$ws.Range("G7").Value = 24
$ws.Range("H7").Value = 68

# Row 9: Estados Unidos -> Estados Unidos
$ws.Range("B9").Value = 18769
$ws.Range("C9").Value = 4980
$ws.Range("E9").Value = 18409
$ws.Range("G9").Value = 28
$ws.Range("H9").Value = 235

# Row 15: Austria -> Austria
$ws.Range("B15").Value = 2649
$ws.Range("C15").Value = 470
$ws.Range("E15").Value = 2634

# Row 17: Noruega -> Noruega
$ws.Range("B17").Value = 1949
$ws.Range("C17").Value = 159
$ws.Range("E17").Value = 1941

# Row 23: Canada -> Canada
$ws.Range("B23").Value = 951
$ws.Range("C23").Value = 78
$ws.Range("E23").Value = 928

# Row 43: Eslovenia -> Eslovenia
$ws.Range("F43").Value = 9

# Row 46: Barein -> Barein
$ws.Range("B46").Value = 298
$ws.Range("C46").Value = 19
$ws.Range("E46").Value = 172

# Row 67: Bulgaria -> Croacia
$ws.Range("A67").Value = "Croacia"
$ws.Range("B67").Value = 130
$ws.Range("C67").Value = 20
$ws.Range("D67").Value = 5
$ws.Range("E67").Value = 124
$ws.Range("H67").Value = 1

# Row 68: Croacia -> Argentina
$ws.Range("A68").Value = "Argentina"
$ws.Range("C68").Value = 0
$ws.Range("D68").Value = 3
$ws.Range("H68").Value = 3

# Row 69: Argentina -> Bulgaria
$ws.Range("A69").Value = "Bulgaria"
$ws.Range("B69").Value = 127
$ws.Range("C69").Value = 20
$ws.Range("D69").Value = 1
$ws.Range("E69").Value = 123

# Row 121: Puerto Rico -> Guam
$ws.Range("A121").Value = "Guam"
$ws.Range("C121").Value = 2

# Row 122: Guam -> Puerto Rico
$ws.Range("A122").Value = "Puerto Rico"
$ws.Range("C122").Value = 8

# Row 127: Nigeria -> Guatemala
$ws.Range("A127").Value = "Guatemala"
$ws.Range("C127").Value = 3
$ws.Range("D127").Value = 0
$ws.Range("H127").Value = 1

# Row 128: Guatemala -> Nigeria
$ws.Range("A128").Value = "Nigeria"
$ws.Range("C128").Value = 0
$ws.Range("D128").Value = 1
$ws.Range("H128").Value = 0

# Row 129: Monaco -> Polinesia Francesa
$ws.Range("A129").Value = "Polinesia Francesa"
$ws.Range("C129").Value = 5

# Row 130: Polinesia Francesa -> Monaco
$ws.Range("A130").Value = "Monaco"
$ws.Range("C130").Value = 1

# Row 132: Trinidad yTobago -> Togo
$ws.Range("A132").Value = "Togo"
$ws.Range("C132").Value = 8

# Row 133: Togo -> Trinidad yTobago
$ws.Range("A133").Value = "Trinidad yTobago"
$ws.Range("C133").Value = 0

# Row 138: Guinea Ecuatorial -> Tanzania
$ws.Range("A138").Value = "Tanzania"

# Row 139: Mongolia -> Kirguistan
$ws.Range("A139").Value = "Kirguistan"
$ws.Range("C139").Value = 3

# Row 140: Tanzania -> Mayotte
$ws.Range("A140").Value = "Mayotte"
$ws.Range("C140").Value = 2

# Row 141: Mayotte -> Mongolia
$ws.Range("A141").Value = "Mongolia"
$ws.Range("C141").Value = 0

# Row 142: Kirguistan -> Guinea Ecuatorial
$ws.Range("A142").Value = "Guinea Ecuatorial"
$ws.Range("C142").Value = 0

# Row 145: San Martin (Parte Francesa) -> Bahamas
$ws.Range("A145").Value = "Bahamas"

# Row 146: Bahamas -> San Martin (Parte Francesa)
$ws.Range("A146").Value = "San Martin (Parte Francesa)"

# Row 148: San Bartolome -> Madagascar
$ws.Range("A148").Value = "Madagascar"
$ws.Range("C148").Value = 3

# Row 149: Islas Virgenes de los Estados Unidos -> San Bartolome
$ws.Range("A149").Value = "San Bartolome"

# Row 150: Namibia -> Congo
$ws.Range("A150").Value = "Congo"

# Row 151: Madagascar -> Islas Virgenes de los Estados Unidos
$ws.Range("A151").Value = "Islas Virgenes de los Estados Unidos"
$ws.Range("C151").Value = 0

# Row 152: Republica de Africa Central -> Namibia
$ws.Range("A152").Value = "Namibia"
$ws.Range("C152").Value = 0

# Row 153: Congo -> Republica de Africa Central
$ws.Range("A153").Value = "Republica de Africa Central"
$ws.Range("C153").Value = 2

# Row 157: Liberia -> Guinea
$ws.Range("A157").Value = "Guinea"
$ws.Range("C157").Value = 1

# Row 159: Haiti -> Zambia
$ws.Range("A159").Value = "Zambia"
$ws.Range("C159").Value = 0

# Row 160: Santa Lucia -> Butan
$ws.Range("A160").Value = "Butan"
$ws.Range("C160").Value = 1

# Row 162: Groenlandia -> Haiti
$ws.Range("A162").Value = "Haiti"
$ws.Range("C162").Value = 2

# Row 163: Benin -> Nicaragua
$ws.Range("A163").Value = "Nicaragua"
$ws.Range("C163").Value = 1

# Row 164: Zambia -> Groenlandia
$ws.Range("A164").Value = "Groenlandia"

# Row 165: Mauritania -> Nueva Caledonia
$ws.Range("A165").Value = "Nueva Caledonia"

# Row 166: Guinea -> Mauritania
$ws.Range("A166").Value = "Mauritania"
$ws.Range("C166").Value = 0

# Row 167: Nueva Caledonia -> Santa Lucia
$ws.Range("A167").Value = "Santa Lucia"

# Row 168: Butan -> Liberia
$ws.Range("A168").Value = "Liberia"
$ws.Range("C168").Value = 0

# Row 169: Nicaragua -> Benin
$ws.Range("A169").Value = "Benin"
$ws.Range("C169").Value = 0

# Row 171: Gambia -> Suazilandia
$ws.Range("A171").Value = "Suazilandia"

# Row 172: Suazilandia -> Angola
$ws.Range("A172").Value = "Angola"
$ws.Range("C172").Value = 1

# Row 173: San Vicente y las Granadinas -> Zimbabue
$ws.Range("A173").Value = "Zimbabue"
$ws.Range("C173").Value = 1

# Row 174: Papua Nueva Guinea -> Antigua y Barbuda
$ws.Range("A174").Value = "Antigua y Barbuda"
$ws.Range("C174").Value = 0

# Row 175: Republica de Yibuti -> San Vicente y las Granadinas
$ws.Range("A175").Value = "San Vicente y las Granadinas"

# Row 176: Santa Sede -> Fiyi
$ws.Range("A176").Value = "Fiyi"

# Row 177: San Martin (Parte Holandesa) -> Santa Sede
$ws.Range("A177").Value = "Santa Sede"

# Row 178: Somalia -> Republica del Chad
$ws.Range("A178").Value = "Republica del Chad"

# Row 179: Zimbabue -> Somalia
$ws.Range("A179").Value = "Somalia"
$ws.Range("C179").Value = 0

# Row 180: El Salvador -> Gambia
$ws.Range("A180").Value = "Gambia"

# Row 181: Republica del Chad -> Niger
$ws.Range("A181").Value = "Niger"

# Row 182: Niger -> Papua Nueva Guinea
$ws.Range("A182").Value = "Papua Nueva Guinea"
$ws.Range("C182").Value = 1

# Row 184: Fiyi -> Cabo Verde
$ws.Range("A184").Value = "Cabo Verde"
$ws.Range("C184").Value = 1

# Row 185: Cabo Verde -> Republica de Yibuti
$ws.Range("A185").Value = "Republica de Yibuti"
$ws.Range("C185").Value = 0

# Row 186: Angola -> San Martin (Parte Holandesa)
$ws.Range("A186").Value = "San Martin (Parte Holandesa)"
$ws.Range("C186").Value = 0

# Row 187: Antigua y Barbuda -> El Salvador
$ws.Range("A187").Value = "El Salvador"

